# Auto-generated Excel COM-interop script
# Applies the 'Added Mounted Tire Processing Pipeline' recomputation
# to Step1_Data (segments 1, 5, 10), and the downstream Step2_Sj
# cumulative-sum sheet and Step3_DataPts_* summary sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Cells.Item(2, 4).Value = 0.002006765722526697  # D2: 0 -> 0.002006765722526697
$ws.Cells.Item(2, 5).Value = 0.2126047528472678  # E2: 0.2127069668866213 -> 0.2126047528472678
$ws.Cells.Item(2, 6).Value = 0.09858344942258712  # F2: 0.09863084541182979 -> 0.09858344942258712
$ws.Cells.Item(2, 7).Value = 0.02239882578945404  # G2: 0.02240959447844173 -> 0.02239882578945404
$ws.Cells.Item(2, 9).Value = 0.0886589680411756  # I2: 0.08870159263506186 -> 0.0886589680411756
$ws.Cells.Item(2, 10).Value = 0.002833018654792688  # J2: 0.002834380685868628 -> 0.002833018654792688
$ws.Cells.Item(2, 11).Value = 0.01047322036895299  # K2: 0.01047825558168749 -> 0.01047322036895299
$ws.Cells.Item(2, 12).Value = 0.004993591018604057  # L2: 0.004995991788587114 -> 0.004993591018604057
$ws.Cells.Item(2, 13).Value = 0.1616269432058169  # M2: 0.1617046486310822 -> 0.1616269432058169
$ws.Cells.Item(2, 14).Value = 0.06101121892154572  # N2: 0.06104055130028296 -> 0.06101121892154572
$ws.Cells.Item(2, 15).Value = 0.004612244023495473  # O2: 0.004614461453189832 -> 0.004612244023495473
$ws.Cells.Item(2, 16).Value = 0.002017067364935771  # P2: 0.00201803711090926 -> 0.002017067364935771
$ws.Cells.Item(2, 20).Value = 0.06754673797541313  # T2: 0.06757921243725445 -> 0.06754673797541313
$ws.Cells.Item(2, 21).Value = 0.04097535745988592  # U2: 0.04099505719257693 -> 0.04097535745988592
$ws.Cells.Item(2, 22).Value = 0.02273830750405345  # V2: 0.02274923940574862 -> 0.02273830750405345
$ws.Cells.Item(2, 23).Value = 0.002626035230761136  # W2: 0.002627297750365359 -> 0.002626035230761136
$ws.Cells.Item(2, 24).Value = 0.009582816705986956  # X2: 0.009587423839133254 -> 0.009582816705986956
$ws.Cells.Item(2, 25).Value = 0.00139260357109053  # Y2: 0.001393273093452157 -> 0.00139260357109053
$ws.Cells.Item(2, 26).Value = 0.005761801011625369  # Z2: 0.005764571114115846 -> 0.005761801011625369
$ws.Cells.Item(2, 27).Value = 0.01274592809790704  # AA2: 0.01275205596089575 -> 0.01274592809790704
$ws.Cells.Item(2, 28).Value = 0.01482724598123254  # AB2: 0.01483437447993228 -> 0.01482724598123254
$ws.Cells.Item(2, 29).Value = 0.0006943343072662767  # AC2: 0.0006946681225420735 -> 0.0006943343072662767
$ws.Cells.Item(2, 30).Value = 0.05960118523007809  # AD2: 0.05962983970657044 -> 0.05960118523007809
$ws.Cells.Item(2, 31).Value = 0.04104178214460991  # AE2: 0.04106151381231295 -> 0.04104178214460991
$ws.Cells.Item(2, 32).Value = 0.01659738254140204  # AF2: 0.01660536206909148 -> 0.01659738254140204
$ws.Cells.Item(2, 33).Value = 0.02584386263194292  # AG2: 0.02585628759214127 -> 0.02584386263194292
$ws.Cells.Item(2, 34).Value = 0.006204554225589697  # AH2: 0.006207537190651822 -> 0.006204554225589697
$ws.Cells.Item(2, 36).Value = 0  # AJ2: 0.001526960269653014 -> 0
$ws.Cells.Item(3, 6).Value = 0.003589283571538826  # F3: 0.003589283571538827 -> 0.003589283571538826
$ws.Cells.Item(3, 10).Value = 0.004894497176005509  # J3: 0.004894497176005506 -> 0.004894497176005509
$ws.Cells.Item(3, 13).Value = 0.2096146787542454  # M3: 0.2096146787542453 -> 0.2096146787542454
$ws.Cells.Item(3, 14).Value = 0.01234425148966899  # N3: 0.01234425148966898 -> 0.01234425148966899
$ws.Cells.Item(3, 15).Value = 0.04440259957954139  # O3: 0.04440259957954138 -> 0.04440259957954139
$ws.Cells.Item(3, 20).Value = 0.06235320421401433  # T3: 0.06235320421401432 -> 0.06235320421401433
$ws.Cells.Item(3, 22).Value = 0.02684863933364153  # V3: 0.02684863933364152 -> 0.02684863933364153
$ws.Cells.Item(3, 32).Value = 0.0350982311658276  # AF3: 0.03509823116582758 -> 0.0350982311658276
$ws.Cells.Item(6, 5).Value = 0.1874019429862166  # E6: 0.1867706865320108 -> 0.1874019429862166
$ws.Cells.Item(6, 6).Value = 0.09581990994206543  # F6: 0.09549714414984097 -> 0.09581990994206543
$ws.Cells.Item(6, 7).Value = 0.003463561979926867  # G6: 0.003451895100600371 -> 0.003463561979926867
$ws.Cells.Item(6, 9).Value = 0.04465733011515495  # I6: 0.04450690356453555 -> 0.04465733011515495
$ws.Cells.Item(6, 11).Value = 0.001010143745888983  # K6: 0.001006741114362826 -> 0.001010143745888983
$ws.Cells.Item(6, 12).Value = 0.05639352785071134  # L6: 0.05620356835582065 -> 0.05639352785071134
$ws.Cells.Item(6, 13).Value = 0.2207006329505972  # M6: 0.2199572111014031 -> 0.2207006329505972
$ws.Cells.Item(6, 14).Value = 0.002967728020386794  # N6: 0.002957731339256682 -> 0.002967728020386794
$ws.Cells.Item(6, 15).Value = 0.00245099772258602  # O6: 0.002442741628188215 -> 0.00245099772258602
$ws.Cells.Item(6, 20).Value = 0.1384483953574798  # T6: 0.1379820371023241 -> 0.1384483953574798
$ws.Cells.Item(6, 21).Value = 0.0310403807790035  # U6: 0.03093582241426343 -> 0.0310403807790035
$ws.Cells.Item(6, 22).Value = 0.009025621767844444  # V6: 0.008995219297606495 -> 0.009025621767844444
$ws.Cells.Item(6, 24).Value = 0.0367436089338074  # X6: 0.03661983945777809 -> 0.0367436089338074
$ws.Cells.Item(6, 26).Value = 0.005083074860482158  # Z6: 0.005065952712430966 -> 0.005083074860482158
$ws.Cells.Item(6, 27).Value = 0.01051292805777735  # AA6: 0.01047751565178345 -> 0.01051292805777735
$ws.Cells.Item(6, 28).Value = 0.008484581257039542  # AB6: 0.00845600126157955 -> 0.008484581257039542
$ws.Cells.Item(6, 29).Value = 0.001626026957728258  # AC6: 0.001620549746577602 -> 0.001626026957728258
$ws.Cells.Item(6, 30).Value = 0.07088216159092038  # AD6: 0.07064339767375213 -> 0.07088216159092038
$ws.Cells.Item(6, 31).Value = 0.03988679218028005  # AE6: 0.03975243500873663 -> 0.03988679218028005
$ws.Cells.Item(6, 32).Value = 0.01929945320581377  # AF6: 0.0192344437176266 -> 0.01929945320581377
$ws.Cells.Item(6, 33).Value = 0.007448104050854063  # AG6: 0.007423015390198904 -> 0.007448104050854063
$ws.Cells.Item(6, 34).Value = 0.006653095687435127  # AH6: 0.006630684982795576 -> 0.006653095687435127
$ws.Cells.Item(6, 36).Value = 0  # AJ6: 0.003368462696527795 -> 0
$ws.Cells.Item(11, 4).Value = 0.003450276134081306  # D11: 0 -> 0.003450276134081306
$ws.Cells.Item(11, 5).Value = 0.1892852887046901  # E11: 0.1874767773268845 -> 0.1892852887046901
$ws.Cells.Item(11, 6).Value = 0.1274769847492393  # F11: 0.1262590159419168 -> 0.1274769847492393
$ws.Cells.Item(11, 7).Value = 0.07688609821689  # G11: 0.07615149604906249 -> 0.07688609821689
$ws.Cells.Item(11, 8).Value = 0.01132398717824051  # H11: 0.0112157930349232 -> 0.01132398717824051
$ws.Cells.Item(11, 9).Value = 0.00834318677091108  # I11: 0.008263472450238804 -> 0.00834318677091108
$ws.Cells.Item(11, 10).Value = 0.003333264945660729  # J11: 0.003301417528353657 -> 0.003333264945660729
$ws.Cells.Item(11, 11).Value = 0.007756763382536089  # K11: 0.007682651997925804 -> 0.007756763382536089
$ws.Cells.Item(11, 12).Value = 0.1240457083001149  # L11: 0.1228605233532873 -> 0.1240457083001149
$ws.Cells.Item(11, 13).Value = 0.01186735533165446  # M11: 0.01175396962012543 -> 0.01186735533165446
$ws.Cells.Item(11, 14).Value = 0.06153745616132051  # N11: 0.06094950138474729 -> 0.06153745616132051
$ws.Cells.Item(11, 15).Value = 0.07738202044372576  # O11: 0.0766426800260539 -> 0.07738202044372576
$ws.Cells.Item(11, 16).Value = 0.01017937668065492  # P11: 0.01008211862815695 -> 0.01017937668065492
$ws.Cells.Item(11, 20).Value = 0.05221404992449518  # T11: 0.05171517489825327 -> 0.05221404992449518
$ws.Cells.Item(11, 21).Value = 0.05777622681757354  # U11: 0.0572242084104311 -> 0.05777622681757354
$ws.Cells.Item(11, 22).Value = 0.01372811795363715  # V11: 0.01359695373223934 -> 0.01372811795363715
$ws.Cells.Item(11, 24).Value = 0.003794274782814895  # X11: 0.00375802268334008 -> 0.003794274782814895
$ws.Cells.Item(11, 25).Value = 0.005388235854459357  # Y11: 0.005336754379508112 -> 0.005388235854459357
$ws.Cells.Item(11, 26).Value = 0.0009283449159483397  # Z11: 0.000919475117589948 -> 0.0009283449159483397
$ws.Cells.Item(11, 27).Value = 0.01107327838962849  # AA11: 0.01096747962367951 -> 0.01107327838962849
$ws.Cells.Item(11, 29).Value = 0.006320198646791489  # AC11: 0.006259812806766983 -> 0.006320198646791489
$ws.Cells.Item(11, 30).Value = 0.0634266902037722  # AD11: 0.06282068488938651 -> 0.0634266902037722
$ws.Cells.Item(11, 31).Value = 0.06239399775178545  # AE11: 0.06179785921607003 -> 0.06239399775178545
$ws.Cells.Item(11, 32).Value = 0.004680248870661387  # AF11: 0.004635531801567031 -> 0.004680248870661387
$ws.Cells.Item(11, 35).Value = 0.005408568888712818  # AI11: 0.005356893143387767 -> 0.005408568888712818
$ws.Cells.Item(11, 36).Value = 0  # AJ11: 0.01297173195610424 -> 0

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Cells.Item(2, 4).Value = 0.002006765722526697  # D2: 0 -> 0.002006765722526697
$ws.Cells.Item(2, 5).Value = 0.2146115185697945  # E2: 0.2127069668866213 -> 0.2146115185697945
$ws.Cells.Item(2, 6).Value = 0.3131949679923817  # F2: 0.3113378122984511 -> 0.3131949679923817
$ws.Cells.Item(2, 7).Value = 0.3355937937818357  # G2: 0.3337474067768929 -> 0.3355937937818357
$ws.Cells.Item(2, 8).Value = 0.3355937937818357  # H2: 0.3337474067768929 -> 0.3355937937818357
$ws.Cells.Item(2, 9).Value = 0.4242527618230113  # I2: 0.4224489994119547 -> 0.4242527618230113
$ws.Cells.Item(2, 10).Value = 0.427085780477804  # J2: 0.4252833800978233 -> 0.427085780477804
$ws.Cells.Item(2, 11).Value = 0.437559000846757  # K2: 0.4357616356795108 -> 0.437559000846757
$ws.Cells.Item(2, 12).Value = 0.4425525918653611  # L2: 0.4407576274680979 -> 0.4425525918653611
$ws.Cells.Item(2, 13).Value = 0.604179535071178  # M2: 0.6024622760991801 -> 0.604179535071178
$ws.Cells.Item(2, 14).Value = 0.6651907539927238  # N2: 0.6635028273994631 -> 0.6651907539927238
$ws.Cells.Item(2, 15).Value = 0.6698029980162192  # O2: 0.6681172888526529 -> 0.6698029980162192
$ws.Cells.Item(2, 16).Value = 0.671820065381155  # P2: 0.6701353259635622 -> 0.671820065381155
$ws.Cells.Item(2, 17).Value = 0.671820065381155  # Q2: 0.6701353259635622 -> 0.671820065381155
$ws.Cells.Item(2, 18).Value = 0.671820065381155  # R2: 0.6701353259635622 -> 0.671820065381155
$ws.Cells.Item(2, 19).Value = 0.671820065381155  # S2: 0.6701353259635622 -> 0.671820065381155
$ws.Cells.Item(2, 20).Value = 0.7393668033565681  # T2: 0.7377145384008166 -> 0.7393668033565681
$ws.Cells.Item(2, 21).Value = 0.780342160816454  # U2: 0.7787095955933935 -> 0.780342160816454
$ws.Cells.Item(2, 22).Value = 0.8030804683205075  # V2: 0.8014588349991421 -> 0.8030804683205075
$ws.Cells.Item(2, 23).Value = 0.8057065035512686  # W2: 0.8040861327495075 -> 0.8057065035512686
$ws.Cells.Item(2, 24).Value = 0.8152893202572555  # X2: 0.8136735565886407 -> 0.8152893202572555
$ws.Cells.Item(2, 25).Value = 0.816681923828346  # Y2: 0.8150668296820929 -> 0.816681923828346
$ws.Cells.Item(2, 26).Value = 0.8224437248399714  # Z2: 0.8208314007962088 -> 0.8224437248399714
$ws.Cells.Item(2, 27).Value = 0.8351896529378785  # AA2: 0.8335834567571045 -> 0.8351896529378785
$ws.Cells.Item(2, 28).Value = 0.850016898919111  # AB2: 0.8484178312370368 -> 0.850016898919111
$ws.Cells.Item(2, 29).Value = 0.8507112332263773  # AC2: 0.8491124993595789 -> 0.8507112332263773
$ws.Cells.Item(2, 30).Value = 0.9103124184564554  # AD2: 0.9087423390661493 -> 0.9103124184564554
$ws.Cells.Item(2, 31).Value = 0.9513542006010653  # AE2: 0.9498038528784623 -> 0.9513542006010653
$ws.Cells.Item(2, 32).Value = 0.9679515831424674  # AF2: 0.9664092149475537 -> 0.9679515831424674
$ws.Cells.Item(2, 33).Value = 0.9937954457744104  # AG2: 0.992265502539695 -> 0.9937954457744104
$ws.Cells.Item(2, 34).Value = 1  # AH2: 0.9984730397303468 -> 1
$ws.Cells.Item(2, 35).Value = 1  # AI2: 0.9984730397303468 -> 1
$ws.Cells.Item(2, 36).Value = 1  # AJ2: 0.9999999999999999 -> 1
$ws.Cells.Item(3, 8).Value = 0.2499601219552052  # H3: 0.2499601219552051 -> 0.2499601219552052
$ws.Cells.Item(3, 31).Value = 0.9012657964565258  # AE3: 0.9012657964565257 -> 0.9012657964565258
$ws.Cells.Item(3, 32).Value = 0.9363640276223534  # AF3: 0.9363640276223533 -> 0.9363640276223534
$ws.Cells.Item(3, 33).Value = 0.9851456292440509  # AG3: 0.9851456292440508 -> 0.9851456292440509
$ws.Cells.Item(6, 5).Value = 0.1874019429862166  # E6: 0.1867706865320108 -> 0.1874019429862166
$ws.Cells.Item(6, 6).Value = 0.2832218529282821  # F6: 0.2822678306818517 -> 0.2832218529282821
$ws.Cells.Item(6, 7).Value = 0.2866854149082089  # G6: 0.2857197257824521 -> 0.2866854149082089
$ws.Cells.Item(6, 8).Value = 0.2866854149082089  # H6: 0.2857197257824521 -> 0.2866854149082089
$ws.Cells.Item(6, 9).Value = 0.3313427450233639  # I6: 0.3302266293469877 -> 0.3313427450233639
$ws.Cells.Item(6, 10).Value = 0.3313427450233639  # J6: 0.3302266293469877 -> 0.3313427450233639
$ws.Cells.Item(6, 11).Value = 0.3323528887692528  # K6: 0.3312333704613505 -> 0.3323528887692528
$ws.Cells.Item(6, 12).Value = 0.3887464166199642  # L6: 0.3874369388171711 -> 0.3887464166199642
$ws.Cells.Item(6, 13).Value = 0.6094470495705614  # M6: 0.6073941499185742 -> 0.6094470495705614
$ws.Cells.Item(6, 14).Value = 0.6124147775909482  # N6: 0.6103518812578309 -> 0.6124147775909482
$ws.Cells.Item(6, 15).Value = 0.6148657753135343  # O6: 0.6127946228860192 -> 0.6148657753135343
$ws.Cells.Item(6, 16).Value = 0.6148657753135343  # P6: 0.6127946228860192 -> 0.6148657753135343
$ws.Cells.Item(6, 17).Value = 0.6148657753135343  # Q6: 0.6127946228860192 -> 0.6148657753135343
$ws.Cells.Item(6, 18).Value = 0.6148657753135343  # R6: 0.6127946228860192 -> 0.6148657753135343
$ws.Cells.Item(6, 19).Value = 0.6148657753135343  # S6: 0.6127946228860192 -> 0.6148657753135343
$ws.Cells.Item(6, 20).Value = 0.753314170671014  # T6: 0.7507766599883432 -> 0.753314170671014
$ws.Cells.Item(6, 21).Value = 0.7843545514500175  # U6: 0.7817124824026066 -> 0.7843545514500175
$ws.Cells.Item(6, 22).Value = 0.793380173217862  # V6: 0.7907077017002131 -> 0.793380173217862
$ws.Cells.Item(6, 23).Value = 0.793380173217862  # W6: 0.7907077017002131 -> 0.793380173217862
$ws.Cells.Item(6, 24).Value = 0.8301237821516694  # X6: 0.8273275411579912 -> 0.8301237821516694
$ws.Cells.Item(6, 25).Value = 0.8301237821516694  # Y6: 0.8273275411579912 -> 0.8301237821516694
$ws.Cells.Item(6, 26).Value = 0.8352068570121516  # Z6: 0.8323934938704222 -> 0.8352068570121516
$ws.Cells.Item(6, 27).Value = 0.845719785069929  # AA6: 0.8428710095222056 -> 0.845719785069929
$ws.Cells.Item(6, 28).Value = 0.8542043663269685  # AB6: 0.8513270107837851 -> 0.8542043663269685
$ws.Cells.Item(6, 29).Value = 0.8558303932846968  # AC6: 0.8529475605303627 -> 0.8558303932846968
$ws.Cells.Item(6, 30).Value = 0.9267125548756172  # AD6: 0.9235909582041149 -> 0.9267125548756172
$ws.Cells.Item(6, 31).Value = 0.9665993470558972  # AE6: 0.9633433932128515 -> 0.9665993470558972
$ws.Cells.Item(6, 32).Value = 0.985898800261711  # AF6: 0.9825778369304781 -> 0.985898800261711
$ws.Cells.Item(6, 33).Value = 0.993346904312565  # AG6: 0.9900008523206769 -> 0.993346904312565
$ws.Cells.Item(6, 34).Value = 1  # AH6: 0.9966315373034725 -> 1
$ws.Cells.Item(6, 35).Value = 1  # AI6: 0.9966315373034725 -> 1
$ws.Cells.Item(11, 4).Value = 0.003450276134081306  # D11: 0 -> 0.003450276134081306
$ws.Cells.Item(11, 5).Value = 0.1927355648387714  # E11: 0.1874767773268845 -> 0.1927355648387714
$ws.Cells.Item(11, 6).Value = 0.3202125495880107  # F11: 0.3137357932688012 -> 0.3202125495880107
$ws.Cells.Item(11, 7).Value = 0.3970986478049007  # G11: 0.3898872893178638 -> 0.3970986478049007
$ws.Cells.Item(11, 8).Value = 0.4084226349831412  # H11: 0.401103082352787 -> 0.4084226349831412
$ws.Cells.Item(11, 9).Value = 0.4167658217540523  # I11: 0.4093665548030258 -> 0.4167658217540523
$ws.Cells.Item(11, 10).Value = 0.420099086699713  # J11: 0.4126679723313794 -> 0.420099086699713
$ws.Cells.Item(11, 11).Value = 0.4278558500822491  # K11: 0.4203506243293052 -> 0.4278558500822491
$ws.Cells.Item(11, 12).Value = 0.551901558382364  # L11: 0.5432111476825925 -> 0.551901558382364
$ws.Cells.Item(11, 13).Value = 0.5637689137140185  # M11: 0.5549651173027179 -> 0.5637689137140185
$ws.Cells.Item(11, 14).Value = 0.625306369875339  # N11: 0.6159146186874652 -> 0.625306369875339
$ws.Cells.Item(11, 15).Value = 0.7026883903190648  # O11: 0.6925572987135191 -> 0.7026883903190648
$ws.Cells.Item(11, 16).Value = 0.7128677669997198  # P11: 0.7026394173416761 -> 0.7128677669997198
$ws.Cells.Item(11, 17).Value = 0.7128677669997198  # Q11: 0.7026394173416761 -> 0.7128677669997198
$ws.Cells.Item(11, 18).Value = 0.7128677669997198  # R11: 0.7026394173416761 -> 0.7128677669997198
$ws.Cells.Item(11, 19).Value = 0.7128677669997198  # S11: 0.7026394173416761 -> 0.7128677669997198
$ws.Cells.Item(11, 20).Value = 0.765081816924215  # T11: 0.7543545922399294 -> 0.765081816924215
$ws.Cells.Item(11, 21).Value = 0.8228580437417885  # U11: 0.8115788006503605 -> 0.8228580437417885
$ws.Cells.Item(11, 22).Value = 0.8365861616954257  # V11: 0.8251757543825998 -> 0.8365861616954257
$ws.Cells.Item(11, 23).Value = 0.8365861616954257  # W11: 0.8251757543825998 -> 0.8365861616954257
$ws.Cells.Item(11, 24).Value = 0.8403804364782406  # X11: 0.82893377706594 -> 0.8403804364782406
$ws.Cells.Item(11, 25).Value = 0.8457686723326999  # Y11: 0.834270531445448 -> 0.8457686723326999
$ws.Cells.Item(11, 26).Value = 0.8466970172486482  # Z11: 0.835190006563038 -> 0.8466970172486482
$ws.Cells.Item(11, 27).Value = 0.8577702956382767  # AA11: 0.8461574861867175 -> 0.8577702956382767
$ws.Cells.Item(11, 28).Value = 0.8577702956382767  # AB11: 0.8461574861867175 -> 0.8577702956382767
$ws.Cells.Item(11, 29).Value = 0.8640904942850681  # AC11: 0.8524172989934845 -> 0.8640904942850681
$ws.Cells.Item(11, 30).Value = 0.9275171844888404  # AD11: 0.915237983882871 -> 0.9275171844888404
$ws.Cells.Item(11, 31).Value = 0.9899111822406258  # AE11: 0.9770358430989411 -> 0.9899111822406258
$ws.Cells.Item(11, 32).Value = 0.9945914311112872  # AF11: 0.9816713749005082 -> 0.9945914311112872
$ws.Cells.Item(11, 33).Value = 0.9945914311112872  # AG11: 0.9816713749005082 -> 0.9945914311112872
$ws.Cells.Item(11, 34).Value = 0.9945914311112872  # AH11: 0.9816713749005082 -> 0.9945914311112872
$ws.Cells.Item(11, 35).Value = 1  # AI11: 0.9870282680438959 -> 1

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Cells.Item(2, 6).Value = 0.604179535071178  # F2: 0.6024622760991801 -> 0.604179535071178
$ws.Cells.Item(6, 6).Value = 0.6094470495705614  # F6: 0.6073941499185742 -> 0.6094470495705614
$ws.Cells.Item(11, 6).Value = 0.551901558382364  # F11: 0.5432111476825925 -> 0.551901558382364

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Cells.Item(2, 6).Value = 0.7393668033565681  # F2: 0.7377145384008166 -> 0.7393668033565681
$ws.Cells.Item(6, 6).Value = 0.753314170671014  # F6: 0.7507766599883432 -> 0.753314170671014
$ws.Cells.Item(11, 4).Value = 14  # D11: 15 -> 14
$ws.Cells.Item(11, 6).Value = 0.7026883903190648  # F11: 0.7026394173416761 -> 0.7026883903190648
$ws.Cells.Item(11, 7).Value = 12  # G11: 13 -> 12

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Cells.Item(2, 6).Value = 0.8030804683205075  # F2: 0.8014588349991421 -> 0.8030804683205075
$ws.Cells.Item(6, 6).Value = 0.8301237821516694  # F6: 0.8273275411579912 -> 0.8301237821516694
$ws.Cells.Item(11, 6).Value = 0.8228580437417885  # F11: 0.8115788006503605 -> 0.8228580437417885

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Cells.Item(2, 6).Value = 0.9103124184564554  # F2: 0.9087423390661493 -> 0.9103124184564554
$ws.Cells.Item(3, 6).Value = 0.9012657964565258  # F3: 0.9012657964565257 -> 0.9012657964565258
$ws.Cells.Item(6, 6).Value = 0.9267125548756172  # F6: 0.9235909582041149 -> 0.9267125548756172
$ws.Cells.Item(11, 6).Value = 0.9275171844888404  # F11: 0.915237983882871 -> 0.9275171844888404
